# Update "Förändrad" (C) column dates from 2023-09-14 (45183) to 2023-09-15 (45184)
# for the data rows on the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C6").Value = 45184
